$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.602996254681648
$ws1.Range("C2").Value = 0.5625
$ws1.Range("D2").Value = 0.9269662921348315
$ws1.Range("E2").Value = 0.7001414427157001
$ws1.Range("F2").Value = 0.8206233421750663
$ws1.Range("G2").Value = 0.9044272663387211
$ws1.Range("H2").Value = 0.7682777146544348
$ws1.Range("I2").Value = 495
$ws1.Range("J2").Value = 385
$ws1.Range("K2").Value = 149
$ws1.Range("L2").Value = 39

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.7925531914893617
$ws2.Range("C2").Value = 0.2790262172284644
$ws2.Range("D2").Value = 0.4127423822714681

$ws2.Range("B3").Value = 0.5625
$ws2.Range("C3").Value = 0.9269662921348315
$ws2.Range("D3").Value = 0.7001414427157001

$ws2.Range("B4").Value = 0.602996254681648
$ws2.Range("C4").Value = 0.602996254681648
$ws2.Range("D4").Value = 0.602996254681648
$ws2.Range("E4").Value = 0.602996254681648

$ws2.Range("B5").Value = 0.6775265957446808
$ws2.Range("C5").Value = 0.602996254681648
$ws2.Range("D5").Value = 0.5564419124935841

$ws2.Range("B6").Value = 0.6775265957446808
$ws2.Range("C6").Value = 0.602996254681648
$ws2.Range("D6").Value = 0.5564419124935841

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 149
$ws3.Range("C2").Value = 385
$ws3.Range("B3").Value = 39
$ws3.Range("C3").Value = 495
